# Update lottery results worksheet with the latest draws (rows 305-308).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New draw results to append after the existing last row (304).
$newRows = @(
    @(304, 2, 14, 22, 36, 40, 50, 1, 2),
    @(305, 4, 10, 16, 34, 40, 47, 2, 4),
    @(306, 7, 9, 11, 41, 48, 50, 4, 6),
    @(307, 1, 10, 13, 14, 28, 31, 1, 4)
)

$startRow = 305
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

# Match the updated view: scrolled down a bit, with the newly added rows selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 297
$win.ScrollColumn = 1
$ws.Range("B305:I308").Select() | Out-Null
